$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet Group -> Sheet1
$ws.Name = "Sheet1"

# Insert a new column before column A for the "Table" name, shifting
# GroupID/GroupName/ParticipationLevelID/ActiveFlg from A:D to B:E
$ws.Columns.Item(1).Insert()

# New header for the Query column (set first so shared-string order matches)
$ws.Range("F1").Value = "Query"

# New header + values for the Table column
$ws.Range("A1").Value = "Table"
for ($r = 2; $r -le 11; $r++) {
  $ws.Range("A$r").Value = "[Vol].[tblGroup]"
}

# Build the generated INSERT statement for each data row
for ($r = 2; $r -le 11; $r++) {
  $ws.Range("F$r").Formula = '="INSERT INTO "&A' + $r + '&" ([" &B$1 &"],["&C$1&"],["&D$1&"],["&E$1&"]) VALUES ( ''" & B' + $r + ' & "'',''" & C' + $r + ' & "'',''" & D' + $r + ' & "'' ,''" & E' + $r + ' & "'')"'
  $ws.Range("F$r").WrapText = $false
}

# widen the new Query column to fit the generated SQL text
$ws.Columns.Item(6).ColumnWidth = 82.6666666666667

$ws.Range("F9").Select()
